$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Register")

# ------------------------------------------------------------------
# 1. Register sheet: C1/D1 pick up the same (yellow-fill) header
#    style that A1/B1 already use.
# ------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2. Add the new "signin" worksheet right after "Register".
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "signin"

# Column widths
$newSheet.Columns.Item(1).ColumnWidth = 13.6
$newSheet.Columns.Item(2).ColumnWidth = 13.13

# ------------------------------------------------------------------
# 3. Populate data (cell-write order chosen to reproduce the exact
#    shared-string insertion order of the original edit).
# ------------------------------------------------------------------
$newSheet.Range("A1").Value = "username"
$newSheet.Range("B1").Value = "password"
$newSheet.Range("C1").Value = "exp"

$newSheet.Range("A2").Value = "cinthiya"
$newSheet.Range("B2").Value = 12345
$newSheet.Range("C2").Value = "invalid"

$newSheet.Range("A3").Value = "cinthiya"
$newSheet.Range("C3").Value = "invalid"

$newSheet.Range("C4").Value = "invalid"

$newSheet.Range("B5").Value = "admin@123"
$newSheet.Range("C5").Value = "invalid"

$newSheet.Range("A6").Value = "cin@123"
$newSheet.Range("B6").Value = "cin@123"
$newSheet.Range("C6").Value = "invalid"

$newSheet.Range("A7").Value = "!@@#$"
$newSheet.Range("B7").Value = "!@@#%"
$newSheet.Range("C7").Value = "invalid"

$newSheet.Range("A8").Value = "cinthiya"
$newSheet.Range("B8").Value = "Dsportal@123"
$newSheet.Range("C8").Value = "invalid"

$newSheet.Range("C9").Value = "valid"
$newSheet.Range("A9").Value = "cinthiyaSDET85"
$newSheet.Range("B9").Value = "Dsportal@123"

# ------------------------------------------------------------------
# 4. Header style (A1:B1) on the new sheet - reuse the yellow-fill
#    style already used for Register!A1:B1.
# ------------------------------------------------------------------
$ws.Range("A1:B1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 5. Hyperlinks (mail-to links, mirroring the convention already
#    used on the Register sheet).
# ------------------------------------------------------------------
$newSheet.Hyperlinks.Add($newSheet.Range("B5"), "mailto:admin@123") | Out-Null
$newSheet.Hyperlinks.Add($newSheet.Range("A6"), "mailto:cin@123") | Out-Null
$newSheet.Hyperlinks.Add($newSheet.Range("B6"), "mailto:cin@123") | Out-Null

# Re-apply the Hyperlink cell style (reusing the existing style
# already present in the workbook from Register!C6) after Add(),
# since Hyperlinks.Add() stamps its own formatting on the cell.
$ws.Range("C6").Copy()
$newSheet.Range("B5").PasteSpecial(-4122)
$ws.Range("C6").Copy()
$newSheet.Range("A6").PasteSpecial(-4122)
$ws.Range("C6").Copy()
$newSheet.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 6. Selection / active sheet bookkeeping, matching the target file:
#    Register stays the visible/active tab (selection at G6) while
#    "signin" remembers its own selection (C9).
# ------------------------------------------------------------------
$newSheet.Range("C9").Select()
$ws.Activate()
$ws.Range("G6").Select()
